# Auto-generated edit script updating cryptos.xlsx price/volume data
# per commit "Updated cryptos list on Tue May 30 11:06:58 UTC 2023 with GitHub Actions"
# All D-column (Price) values are forced to Text to match source data,
# which was written as text (t="inlineStr") rather than numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.913.61"
$c.ClearFormats()
$ws.Range("E2").Value = "  -0.35%  "
# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.908.21"
$c.ClearFormats()
$ws.Range("E3").Value = "  -0.20%  "
# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9992"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.76%  "
# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "313.64"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.68%  "
# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.53%  "
# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5018"
$c.ClearFormats()
$ws.Range("E7").Value = "  +4.33%  "
# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3816"
$c.ClearFormats()
$ws.Range("E8").Value = "  +0.20%  "
# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07315"
$c.ClearFormats()
$ws.Range("E9").Value = "  -0.50%  "
# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9109"
$c.ClearFormats()
$ws.Range("E10").Value = "  -2.32%  "
# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "21.18"
$c.ClearFormats()
$ws.Range("E11").Value = "  +1.88%  "
# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07677"
$c.ClearFormats()
$ws.Range("E12").Value = "  -1.63%  "
# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.920.71"
$c.ClearFormats()
$ws.Range("E13").Value = "  +0.50%  "
# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.499"
$c.ClearFormats()
$ws.Range("E14").Value = "  +0.04%  "
# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "92.72"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.98%  "
# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E16").Value = "  -0.59%  "
# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008750"
$c.ClearFormats()
$ws.Range("E17").Value = "  -1.22%  "
# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.57%  "
# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "27.938.83"
$c.ClearFormats()
$ws.Range("E19").Value = "  -0.36%  "
# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.69"
$c.ClearFormats()
$ws.Range("E20").Value = "  -0.29%  "
# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.188"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.47%  "
# Row 22
$ws.Range("B22").Value = "Cosmos"
$ws.Range("C22").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.85"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.50%  "
# Row 23
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.613"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.24%  "
# Row 24
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "153.28"
$c.ClearFormats()
$ws.Range("E24").Value = "  -1.91%  "
# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.849"
$c.ClearFormats()
$ws.Range("E25").Value = "  -2.96%  "
# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.209"
$c.ClearFormats()
$ws.Range("E26").Value = "  +4.10%  "
# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.42"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.34%  "
# Row 28
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "115.48"
$c.ClearFormats()
$ws.Range("E28").Value = "  -1.19%  "
# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "4.926"
$c.ClearFormats()
$ws.Range("E29").Value = "  -0.62%  "
# Row 30
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.09042"
$c.ClearFormats()
$ws.Range("E30").Value = "  +1.06%  "
# Row 31
$ws.Range("B31").Value = "HuobiToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.220"
$c.ClearFormats()
$ws.Range("E31").Value = "  -2.61%  "
# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.896"
$c.ClearFormats()
$ws.Range("E32").Value = "  +4.96%  "
# Row 33
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.233"
$c.ClearFormats()
$ws.Range("E33").Value = "  -1.79%  "
# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7772"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.18%  "
# Row 35
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.02085"
$c.ClearFormats()
$ws.Range("E35").Value = "  +1.66%  "
# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.581"
$c.ClearFormats()
$ws.Range("E36").Value = "  -1.15%  "
# Row 37
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.072"
$c.ClearFormats()
$ws.Range("E37").Value = "  +2.62%  "
# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.095"
$c.ClearFormats()
$ws.Range("E38").Value = "  -1.20%  "
# Row 39
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.5554"
$c.ClearFormats()
$ws.Range("E39").Value = "  +0.88%  "
# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.05298"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.02%  "
# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.888"
$c.ClearFormats()
$ws.Range("E41").Value = "  -1.77%  "
# Row 42
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "113.46"
$c.ClearFormats()
$ws.Range("E42").Value = "  +4.61%  "
# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "8.508"
$c.ClearFormats()
$ws.Range("E43").Value = "  +0.49%  "
# Row 44
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.1522"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.20%  "
# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4834"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.48%  "
# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "10.59"
$c.ClearFormats()
$ws.Range("E46").Value = "  -0.89%  "
# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.56%  "
# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.639"
$c.ClearFormats()
$ws.Range("E48").Value = "  -0.20%  "
# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "67.63"
$c.ClearFormats()
$ws.Range("E49").Value = "  -0.44%  "
# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06062"
$c.ClearFormats()
$ws.Range("E50").Value = "  -0.34%  "
# Row 51
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.9084"
$c.ClearFormats()
$ws.Range("E51").Value = "  +1.20%  "
